{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// 1. Collapse the split runs \"JS\" + \"3\" + \".js contains:\" into a single\n//    run \"JS3.js contains:\" (text content is unchanged, only run layout).\nconst jsResults = body.search(\"JS3.js contains:\");\nawait context.sync();\nif (jsResults.items.length > 0) {\n  jsResults.items[0].insertText(\"JS3.js contains:\", Word.InsertLocation.replace);\n}\n\n// 2. Collapse the split runs \"Commit \" + \"7\" + \":\" into a single run\n//    \"Commit 7:\". There are two \"Commit 7:\" paragraphs in the document;\n//    the one we need is the bullet directly after \"JS3.js contains:\".\nlet commitIndex = -1;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  if (paragraphs.items[i].text.indexOf(\"JS3.js contains\") !== -1) {\n    commitIndex = i + 1;\n    break;\n  }\n}\nif (commitIndex !== -1) {\n  paragraphs.items[commitIndex].getRange().insertText(\"Commit 7:\", Word.InsertLocation.replace);\n}\n\n// 3. Add the new \"Commit 8\" / \"Commit 9\" bullet points right after the\n//    paragraph describing the Object commit, before the trailing blank\n//    bullet.\nlet objectIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"dot and bracket notation and uses of bracket notation\") !== -1) {\n    objectIndex = i;\n    break;\n  }\n}\nif (objectIndex !== -1) {\n  const target = paragraphs.items[objectIndex];\n  const p1 = target.insertParagraph(\"Commit 8:\", Word.InsertLocation.after);\n  const p2 = p1.insertParagraph(\n    \"Iterating object using for in loop, Object.keys(), computed properties, spread operator in objects.\",\n    Word.InsertLocation.after\n  );\n  const p3 = p2.insertParagraph(\"Commit 9:\", Word.InsertLocation.after);\n  p3.insertParagraph(\n    \"Object destructuring, objects inside array, nested destructuring\",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Collapse the split runs \"JS\" + \"3\" + \".js contains:\" into a single\n#    run \"JS3.js contains:\" (same visible text, just re-written as one run).\n$range1 = $d.Content\n$find1 = $range1.Find\n$found1 = $find1.Execute(\"JS3.js contains:\", $false, $false, $false, $false, $false, $true, 1, $false, \"JS3.js contains:\", 2)\n\n# 2. Collapse the split runs \"Commit \" + \"7\" + \":\" into a single run\n#    \"Commit 7:\".\n$range2 = $d.Content\n$find2 = $range2.Find\n$found2 = $find2.Execute(\"Commit 7:\", $false, $false, $false, $false, $false, $true, 1, $false, \"Commit 7:\", 2)\n\n# 3. Add the new \"Commit 8\" / \"Commit 9\" bullet points right after the\n#    paragraph describing the Object commit, before the trailing blank\n#    bullet.\n$range3 = $d.Content\n$find3 = $range3.Find\n$found3 = $find3.Execute(\"dot and bracket notation and uses of bracket notation \")\n$range3.InsertAfter(\"`rCommit 8:`rIterating object using for in loop, Object.keys(), computed properties, spread operator in objects.`rCommit 9:`rObject destructuring, objects inside array, nested destructuring\")\n"}
